# 26-Mar-2021, midday update: record the day's petty-cash transactions for
# 44279 (continuing), 44280 and 44281, and advance the "current" scroll
# position of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---- 44279 (row 19) continued entries -----------------------------------
$ws.Range("B19").Value   = "Wages Expense"
$ws.Range("D19").Formula = "=60000+240000"

$ws.Range("B20").Value   = "TRANSFER BCA"
$ws.Range("D20").Formula = "=600000+1650000"

$ws.Range("B21").Value = "PLN - Astar 214"
$ws.Range("D21").Value = 103000

$ws.Range("B22").Value   = "A/R"
$ws.Range("C22").Formula = "=2250000"

$ws.Range("B23").Value   = "SALES - cash/retail"
$ws.Range("C23").Formula = "=7112525+520475-2250000"

$ws.Range("B24").Value = "SELISIH - kurang"
$ws.Range("D24").Value = 160000

$ws.Range("B25").Value = "SETOR KE BANK"
$ws.Range("D25").Value = 5000000

# ---- 44280 (row 26) ------------------------------------------------------
$ws.Range("A26").Value = 44280
$ws.Range("A26").NumberFormat = $ws.Range("A19").NumberFormat

$ws.Range("B26").Value   = "Wages Expense"
$ws.Range("D26").Formula = "=60000+260000"

$ws.Range("B27").Value   = "FREIGHT OUT"
$ws.Range("D27").Formula = "=12000"

$ws.Range("B28").Value   = "TRANSFER BCA"
$ws.Range("D28").Formula = "=2396000+43000000"

$ws.Range("B29").Value   = "A/R"
$ws.Range("C29").Formula = "=43000000+26143000"

$ws.Range("B30").Value   = "SALES - cash/retail"
$ws.Range("C30").Formula = "=2405475+33554025-26143000"

$ws.Range("B31").Value = "SETOR KE BANK"
$ws.Range("D31").Value = 33000000

# ---- 44281 (row 32) ------------------------------------------------------
$ws.Range("A32").Value = 44281
$ws.Range("A32").NumberFormat = $ws.Range("A19").NumberFormat

$ws.Range("B32").Value   = "Wages Expense"
$ws.Range("D32").Formula = "=60000"

$ws.Range("B33").Value   = "Materai"
$ws.Range("D33").Formula = "=20000"

$ws.Range("B34").Value   = "BELI lampu"
$ws.Range("D34").Formula = "=35000"

$ws.Range("B35").Value   = "TRANSFER BCA"
$ws.Range("D35").Formula = "=2690000"

$ws.Range("B36").Value   = "LPG"
$ws.Range("D36").Formula = "=145000"

# ---- scroll / selection position ----------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 31
$ws.Range("C52").Select() | Out-Null
